# Update "想去人数" (want-to-go headcount) figures across the workbook.
# Generated output refresh (gh-pages data update).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 54957
$ws.Range("F4").Value = 3117
$ws.Range("F6").Value = 1181
$ws.Range("F8").Value = 860
$ws.Range("F10").Value = 1102
$ws.Range("F11").Value = 1367
$ws.Range("F12").Value = 121
$ws.Range("F14").Value = 220
$ws.Range("F16").Value = 57
$ws.Range("F21").Value = 5413
$ws.Range("F23").Value = 5297
$ws.Range("F24").Value = 9313
$ws.Range("F27").Value = 150
$ws.Range("F28").Value = 240
$ws.Range("F29").Value = 453
$ws.Range("F31").Value = 107
$ws.Range("F33").Value = 288

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 140

# --- Sheet 3: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 594
$ws.Range("F5").Value = 52

# --- Sheet 4: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 594
$ws.Range("F5").Value = 3117
$ws.Range("F6").Value = 52
$ws.Range("F7").Value = 140
$ws.Range("F8").Value = 1181
$ws.Range("F11").Value = 860
$ws.Range("F13").Value = 1102
$ws.Range("F15").Value = 1367
$ws.Range("F17").Value = 121
$ws.Range("F18").Value = 220
$ws.Range("F21").Value = 57
$ws.Range("F26").Value = 5414
$ws.Range("F28").Value = 5298
$ws.Range("F29").Value = 9313
$ws.Range("F33").Value = 150
$ws.Range("F34").Value = 240
$ws.Range("F35").Value = 453
$ws.Range("F47").Value = 288
